# daily auto push: 2026-01-15 13:50 UTC
# A new observation for 2026/01/15 (time 20) was logged. This inserts a new
# row at row 645 (pushing the existing rows 645:686 down to 646:687) and
# fills in the new row's data, matching how the rest of the "log" sheet is
# laid out (date / weekday / hour / ranking).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 645:686 down to 646:687 and make room for the new entry.
$ws.Rows.Item(645).Insert()

# Leading apostrophe keeps "2026/01/15" as literal text instead of having it
# auto-converted into a date serial value; resetting the style afterwards
# drops the resulting "quote prefix" number format so the cell keeps the
# sheet's default (unstyled) look, same as every other data row.
$ws.Range("A645").Value = "'2026/01/15"
$ws.Range("A645").Style = "Normal"

$ws.Range("B645").Value = "木"
$ws.Range("C645").Value = 20
$ws.Range("D645").Value = 201
